# "more fix to tex"
#
# The slide 2 diagram shows the formula "Angular velocity omega_b, bias b_omega".
# The final (baseline-raised) character of that run - the Greek "omega" used as
# the bias subscript - is corrected to a Latin "f" (bias b_f).
#
# The run is nested several group levels deep (Group 22 > Group 17 > Group 4 >
# TextBox 28), so we look the shape up recursively by name instead of hard
# coding GroupItems indices.

function Find-ShapeByName($shapes, $name) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -eq $name) {
            return $shp
        }
        if ($shp.Type -eq 6) {
            $found = Find-ShapeByName $shp.GroupItems $name
            if ($found -ne $null) {
                return $found
            }
        }
    }
    return $null
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

$shp = Find-ShapeByName $s.Shapes "TextBox 28"
$tr = $shp.TextFrame.TextRange

# Locate the trailing omega (U+03C9) in "...bias b<omega> " and change it to "f".
$omega = [char]0x03C9
$text = $tr.Text
$pos0 = $text.LastIndexOf($omega)
$pos1 = $pos0 + 1

$tr.Characters($pos1, 1).Text = "f"
